# complete new trade mechanics
# Update stock prices/values/RSI on the "stocks" sheet following new trades,
# and roll the updated STOCKS total into the "portfolio" sheet.

$wb = $excel.ActiveWorkbook

$stocks = $wb.Worksheets.Item("stocks")

# MSFT row
$stocks.Range("B2").Value = 177.36
$stocks.Range("D2").Value = 177.36
$stocks.Range("E2").Value = 54.94932666944332

# AAPL row
$stocks.Range("B3").Value = 295.7
$stocks.Range("D3").Value = 295.7
$stocks.Range("E3").Value = 53.4112618592753

# GOOGL row
$stocks.Range("B4").Value = 1334.05
$stocks.Range("D4").Value = 1334.05
$stocks.Range("E4").Value = 57.45291685806153

$portfolio = $wb.Worksheets.Item("portfolio")

# Updated STOCKS total value
$portfolio.Range("B3").Value = 1807
